# Final set of survival and size data
# Fill in the remaining "length" (column E) measurements for the BSW
# treatment cups (rows 30-35) on Sheet1, and leave the view positioned
# on the last entered cell, same as the author did when finishing data
# entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E30").Value = 0.915
$ws.Range("E31").Value = 0.811
$ws.Range("E32").Value = 0.825
$ws.Range("E33").Value = 0.793
$ws.Range("E34").Value = "NA"
$ws.Range("E35").Value = 0.786

# Reflect where the author ended up after typing the new values.
$ws.Range("E32").Select() | Out-Null
